# Cancellation.xlsx maintenance edit ("Excel Change Because of Nainsi")
#
# 1. Contract sheet: the sample contractId in A3 is swapped for a new
#    test value.
# 2. Calculate / OverRide sheets: the lingering cell selection left over
#    from editing is moved to A3, and the page is set up as Letter/A4
#    portrait (paperSize 9) like the other sheets in the workbook.

$wb = $excel.ActiveWorkbook

# --- Contract sheet: update the sample contract id -----------------------
$wsContract = $wb.Worksheets.Item("Contract")
$wsContract.Range("A3").Value = "9Z001140A9"

# --- Calculate sheet: reset selection + page setup ------------------------
$wsCalculate = $wb.Worksheets.Item("Calculate")
$wsCalculate.Range("A3").Select() | Out-Null
$wsCalculate.PageSetup.Orientation = 1
$wsCalculate.PageSetup.PaperSize = 9

# --- OverRide sheet: reset selection + page setup -------------------------
$wsOverRide = $wb.Worksheets.Item("OverRide")
$wsOverRide.Range("A3").Select() | Out-Null
$wsOverRide.PageSetup.Orientation = 1
$wsOverRide.PageSetup.PaperSize = 9

# Restore focus to the sheet that was active before we touched the others.
$wsContract.Activate() | Out-Null
